# "also add a time variable." -- fills in the previously-blank DATA IN /
# DATA OUT table (rows 6-10) with the five test-case scenarios and their
# population-projection formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Regular US Data
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 126
$ws.Range("E6").Value = 33310036
$ws.Range("F6").Value = 5

# Row 7: Population Increase - High Birth Rate
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = 126
$ws.Range("E7").Value = 33310036
$ws.Range("F7").Value = 5

# Row 8: Population Increase - High Migration
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 12
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = 33310036
$ws.Range("F8").Value = 5

# Row 9: Population Decrease - High Death Rate
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 126
$ws.Range("E9").Value = 33310036
$ws.Range("F9").Value = 5

# Row 10: Population Low Birth Rate and Low Migration
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 250
$ws.Range("E10").Value = 33310036
$ws.Range("F10").Value = 5

# Row 6 formulas (entered individually, not as part of the shared group)
$ws.Range("G6").Formula = "=(`$B`$2/B6 + `$B`$2/D6-`$B`$2/C6)*F6"
$ws.Range("H6").Formula = "= E6+G6"
$ws.Range("I6").Formula = "=IF(H6>E6,""increase"",""decrease"")"

# Rows 7:10 formulas, filled as one shared-formula block each
$ws.Range("G7:G10").Formula = "=(`$B`$2/B7 + `$B`$2/D7-`$B`$2/C7)*F7"
$ws.Range("H7:H10").Formula = "= E7+G7"
$ws.Range("I7:I10").Formula = "=IF(H7>E7,""increase"",""decrease"")"

# Columns G:H get auto-fit ("bestFit") once they hold real numbers
$ws.Range("G1:H1").EntireColumn.ColumnWidth = 10.83

# The user scrolled down a bit and left the selection on the empty cell
# just below the table (F11) after finishing data entry.
$ws.Range("F11").Select() | Out-Null
